$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036190752139612
$ws.Range("D2").Value = 1.044154630400617
$ws.Range("E2").Value = 1.044600438516409
$ws.Range("F2").Value = 1.05419676900061
$ws.Range("I2").Value = 1.034310528346633
$ws.Range("J2").Value = 1.041300418226947
$ws.Range("K2").Value = 1.046926478096904
$ws.Range("L2").Value = 1.047371032052837
$ws.Range("M2").Value = 1.056940644597213
$ws.Range("N2").Value = 1.017728136148673

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037072593924504
$ws.Range("D3").Value = 1.044942300383873
$ws.Range("E3").Value = 1.04538887952995
$ws.Range("F3").Value = 1.055079629460442
$ws.Range("I3").Value = 1.034415789990778
$ws.Range("J3").Value = 1.041826462379673
$ws.Range("K3").Value = 1.047525389922442
$ws.Range("L3").Value = 1.047970803474062
$ws.Range("M3").Value = 1.057636522199884
$ws.Range("N3").Value = 1.017902013226897

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037643330380918
$ws.Range("D4").Value = 1.045452445399541
$ws.Range("E4").Value = 1.045899603874789
$ws.Range("F4").Value = 1.055651633680546
$ws.Range("I4").Value = 1.034481736725959
$ws.Range("J4").Value = 1.04216632964216
$ws.Range("K4").Value = 1.047912738532236
$ws.Range("L4").Value = 1.048358788607165
$ws.Range("M4").Value = 1.058086896574145
$ws.Range("N4").Value = 1.018014339468699

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037883296654477
$ws.Range("D5").Value = 1.045667021121109
$ws.Range("E5").Value = 1.046114442402931
$ws.Range("F5").Value = 1.055892278297334
$ws.Range("I5").Value = 1.034508941558404
$ws.Range("J5").Value = 1.042309084398089
$ws.Range("K5").Value = 1.048075533610257
$ws.Range("L5").Value = 1.048521870515746
$ws.Range("M5").Value = 1.058276255161184
$ws.Range("N5").Value = 1.018061516948912

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03792358968616
$ws.Range("D6").Value = 1.045703055773332
$ws.Range("E6").Value = 1.04615052231086
$ws.Range("F6").Value = 1.055932693741139
$ws.Range("I6").Value = 1.034513478896466
$ws.Range("J6").Value = 1.042333046155419
$ws.Range("K6").Value = 1.048102864865141
$ws.Range("L6").Value = 1.048549251059068
$ws.Range("M6").Value = 1.058308050499115
$ws.Range("N6").Value = 1.018069435634253

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037646536711211
$ws.Range("D7").Value = 1.045455312137822
$ws.Range("E7").Value = 1.045902474049035
$ws.Range("F7").Value = 1.055654848504086
$ws.Range("I7").Value = 1.034482102279497
$ws.Range("J7").Value = 1.04216823763276
$ws.Range("K7").Value = 1.047914913991667
$ws.Range("L7").Value = 1.048360967823641
$ws.Range("M7").Value = 1.058089426711512
$ws.Range("N7").Value = 1.018014970031821

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036488747762932
$ws.Range("D8").Value = 1.044420728832614
$ws.Range("E8").Value = 1.044866780872099
$ws.Range("F8").Value = 1.054494982676005
$ws.Range("I8").Value = 1.034346549624798
$ws.Range("J8").Value = 1.041478303877035
$ws.Range("K8").Value = 1.04712892110348
$ws.Range("L8").Value = 1.04757374894179
$ws.Range("M8").Value = 1.057175799083517
$ws.Range("N8").Value = 1.017786936515901

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034449603823397
$ws.Range("D9").Value = 1.042601335526127
$ws.Range("E9").Value = 1.043046047617203
$ws.Range("F9").Value = 1.052456854771239
$ws.Range("I9").Value = 1.034091152673125
$ws.Range("J9").Value = 1.040258639249714
$ws.Range("K9").Value = 1.04574252821356
$ws.Range("L9").Value = 1.046185810311175
$ws.Range("M9").Value = 1.055566666462981
$ws.Range("N9").Value = 1.017383724258817

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03309095096386
$ws.Range("D10").Value = 1.04139097897077
$ws.Range("E10").Value = 1.041835211638713
$ws.Range("F10").Value = 1.05110204304083
$ws.Range("I10").Value = 1.033909826492568
$ws.Range("J10").Value = 1.039442975433712
$ws.Range("K10").Value = 1.044817426146402
$ws.Range("L10").Value = 1.045260090425895
$ws.Range("M10").Value = 1.054494537567976
$ws.Range("N10").Value = 1.017114009564658

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032502840965964
$ws.Range("D11").Value = 1.040867512534391
$ws.Range("E11").Value = 1.041311635757129
$ws.Range("F11").Value = 1.050516351096417
$ws.Range("I11").Value = 1.033828698204864
$ws.Range("J11").Value = 1.039089194160192
$ws.Range("K11").Value = 1.044416664711632
$ws.Range("L11").Value = 1.044859159245036
$ws.Range("M11").Value = 1.054030461372824
$ws.Range("N11").Value = 1.016997010660788

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03228442148629
$ws.Range("D12").Value = 1.040673169520989
$ws.Range("E12").Value = 1.041117266837687
$ws.Range("F12").Value = 1.050298943575119
$ws.Range("I12").Value = 1.033798171884199
$ws.Range("J12").Value = 1.038957696135303
$ws.Range("K12").Value = 1.044267777467823
$ws.Range("L12").Value = 1.044710223632896
$ws.Range("M12").Value = 1.053858108753001
$ws.Range("N12").Value = 1.016953520845894

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032331271792353
$ws.Range("D13").Value = 1.040714852398141
$ws.Range("E13").Value = 1.04115895460443
$ws.Range("F13").Value = 1.050345571657518
$ws.Range("I13").Value = 1.033804737592962
$ws.Range("J13").Value = 1.038985906865997
$ws.Range("K13").Value = 1.044299715469084
$ws.Range("L13").Value = 1.044742171344995
$ws.Range("M13").Value = 1.053895077777083
$ws.Range("N13").Value = 1.016962850964267

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032484785708288
$ws.Range("D14").Value = 1.040851446117713
$ws.Range("E14").Value = 1.041295566898212
$ws.Range("F14").Value = 1.050498377164141
$ws.Range("I14").Value = 1.033826182874355
$ws.Range("J14").Value = 1.039078326278507
$ws.Range("K14").Value = 1.044404358180458
$ws.Range("L14").Value = 1.044846848416229
$ws.Range("M14").Value = 1.054016214113647
$ws.Range("N14").Value = 1.016993416413333

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032579374815106
$ws.Range("D15").Value = 1.040935618751807
$ws.Range("E15").Value = 1.041379752930619
$ws.Range("F15").Value = 1.050592544885701
$ws.Range("I15").Value = 1.033839344147399
$ws.Range("J15").Value = 1.039135257317493
$ws.Range("K15").Value = 1.04446882852403
$ws.Range("L15").Value = 1.044911341875487
$ws.Range("M15").Value = 1.054090853690915
$ws.Range("N15").Value = 1.01701224467234

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033129986126308
$ws.Range("D16").Value = 1.04142573304625
$ws.Range("E16").Value = 1.041869975045672
$ws.Range("F16").Value = 1.051140933660922
$ws.Range("I16").Value = 1.033915155736111
$ws.Range("J16").Value = 1.039466442335281
$ws.Range("K16").Value = 1.044844019557689
$ws.Range("L16").Value = 1.04528669716209
$ws.Range("M16").Value = 1.054525340332526
$ws.Range("N16").Value = 1.017121769995534

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033475423401949
$ws.Range("D17").Value = 1.0417333375704
$ws.Range("E17").Value = 1.042177673480624
$ws.Range("F17").Value = 1.051485179155849
$ws.Range("I17").Value = 1.03396201143808
$ws.Range("J17").Value = 1.039674027948001
$ws.Range("K17").Value = 1.045079318110548
$ws.Range("L17").Value = 1.045522124918564
$ws.Range("M17").Value = 1.05479792694172
$ws.Range("N17").Value = 1.017190416246308

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033676929907263
$ws.Range("D18").Value = 1.041912818366376
$ws.Range("E18").Value = 1.042357218522969
$ws.Range("F18").Value = 1.051686063319217
$ws.Range("I18").Value = 1.033989089414531
$ws.Range("J18").Value = 1.039795051780081
$ws.Range("K18").Value = 1.045216545725257
$ws.Range("L18").Value = 1.04565943732967
$ws.Range("M18").Value = 1.054956937696258
$ws.Range("N18").Value = 1.01723043608805

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033745641563171
$ws.Range("D19").Value = 1.041974026852881
$ws.Range("E19").Value = 1.042418450522035
$ws.Range("F19").Value = 1.05175457510968
$ws.Range("I19").Value = 1.033998279506537
$ws.Range("J19").Value = 1.039836308014571
$ws.Range("K19").Value = 1.045263333692428
$ws.Range("L19").Value = 1.045706255810609
$ws.Range("M19").Value = 1.055011158854178
$ws.Range("N19").Value = 1.017244078354

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033438359298352
$ws.Range("D20").Value = 1.041700328268696
$ws.Range("E20").Value = 1.042144653123669
$ws.Range("F20").Value = 1.051448235385881
$ws.Range("I20").Value = 1.033957010340508
$ws.Range("J20").Value = 1.039651761898575
$ws.Range("K20").Value = 1.045054074659821
$ws.Range("L20").Value = 1.045496866629745
$ws.Range("M20").Value = 1.05476867935413
$ws.Range("N20").Value = 1.0171830532555

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032439578826242
$ws.Range("D21").Value = 1.040811220015715
$ws.Range("E21").Value = 1.041255334919479
$ws.Range("F21").Value = 1.05045337574688
$ws.Range("I21").Value = 1.033819878578192
$ws.Range("J21").Value = 1.039051113481675
$ws.Range("K21").Value = 1.044373544233582
$ws.Range("L21").Value = 1.04481602394596
$ws.Range("M21").Value = 1.053980541753635
$ws.Range("N21").Value = 1.016984416508531

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031811784537034
$ws.Range("D22").Value = 1.040252757053672
$ws.Range("E22").Value = 1.040696825273824
$ws.Range("F22").Value = 1.049828705157227
$ws.Range("I22").Value = 1.033731392792918
$ws.Range("J22").Value = 1.038672953989994
$ws.Range("K22").Value = 1.043945514496589
$ws.Range("L22").Value = 1.044387882828125
$ws.Range("M22").Value = 1.053485159193174
$ws.Range("N22").Value = 1.016859345329745

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032144572781129
$ws.Range("D23").Value = 1.040548755696897
$ws.Range("E23").Value = 1.040992840577033
$ws.Range("F23").Value = 1.05015977492228
$ws.Range("I23").Value = 1.033778515240172
$ws.Range("J23").Value = 1.038873471226911
$ws.Range("K23").Value = 1.044172435232021
$ws.Range("L23").Value = 1.044614854558651
$ws.Range("M23").Value = 1.053747756065572
$ws.Range("N23").Value = 1.01692566488

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033455106913425
$ws.Range("D24").Value = 1.041715243565188
$ws.Range("E24").Value = 1.042159573386296
$ws.Range("F24").Value = 1.051464928402489
$ws.Range("I24").Value = 1.033959270900448
$ws.Range("J24").Value = 1.039661823145183
$ws.Range("K24").Value = 1.045065481144521
$ws.Range("L24").Value = 1.045508279789991
$ws.Range("M24").Value = 1.054781895031974
$ws.Range("N24").Value = 1.017186380337046

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034976640775201
$ws.Range("D25").Value = 1.04307124633241
$ws.Range("E25").Value = 1.043516231874161
$ws.Range("F25").Value = 1.052983072877695
$ws.Range("I25").Value = 1.034159132662095
$ws.Range("J25").Value = 1.040574407567709
$ws.Range("K25").Value = 1.046101098135735
$ws.Range("L25").Value = 1.046544707067572
$ws.Range("M25").Value = 1.057175799083517
$ws.Range("N25").Value = 1.017488126020598
